$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 1).Value = "Última actualización: 08:20:43"
$ws1.Cells.Item(3, 1).Value = "Total filas: 93"

$ws1.Cells.Item(68, 1).Value = "08:20:43"
$ws1.Cells.Item(68, 2).Value = "08:21"
$ws1.Cells.Item(68, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(68, 4).Value = 1
$ws1.Cells.Item(68, 5).Value = "LP1912"
$ws1.Cells.Item(69, 1).Value = "06:59:44"
$ws1.Cells.Item(69, 2).Value = "08:22"
$ws1.Cells.Item(69, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(69, 4).Value = 83
$ws1.Cells.Item(69, 5).Value = "LP1912"
$ws1.Cells.Item(70, 1).Value = "06:44:15"
$ws1.Cells.Item(70, 2).Value = "08:23"
$ws1.Cells.Item(70, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(70, 4).Value = 99
$ws1.Cells.Item(70, 5).Value = "LP1912"
$ws1.Cells.Item(71, 1).Value = "06:44:15"
$ws1.Cells.Item(71, 2).Value = "08:23"
$ws1.Cells.Item(71, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(71, 4).Value = 99
$ws1.Cells.Item(71, 5).Value = "LP1912"
$ws1.Cells.Item(72, 1).Value = "07:57:27"
$ws1.Cells.Item(72, 2).Value = "08:26"
$ws1.Cells.Item(72, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(72, 4).Value = 29
$ws1.Cells.Item(72, 5).Value = "LP1912"
$ws1.Cells.Item(73, 1).Value = "06:44:15"
$ws1.Cells.Item(73, 2).Value = "08:27"
$ws1.Cells.Item(73, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(73, 4).Value = 103
$ws1.Cells.Item(73, 5).Value = "LP1912"
$ws1.Cells.Item(74, 1).Value = "07:57:27"
$ws1.Cells.Item(74, 2).Value = "08:33"
$ws1.Cells.Item(74, 3).Value = "10_OLMOS"
$ws1.Cells.Item(74, 4).Value = 36
$ws1.Cells.Item(74, 5).Value = "LP1912"
$ws1.Cells.Item(75, 1).Value = "08:20:43"
$ws1.Cells.Item(75, 2).Value = "08:37"
$ws1.Cells.Item(75, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(75, 4).Value = 17
$ws1.Cells.Item(75, 5).Value = "LP1912"
$ws1.Cells.Item(76, 1).Value = "06:44:15"
$ws1.Cells.Item(76, 2).Value = "08:42"
$ws1.Cells.Item(76, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(76, 4).Value = 118
$ws1.Cells.Item(76, 5).Value = "LP1912"
$ws1.Cells.Item(77, 1).Value = "07:31:43"
$ws1.Cells.Item(77, 2).Value = "08:43"
$ws1.Cells.Item(77, 3).Value = "14_ABASTO"
$ws1.Cells.Item(77, 4).Value = 72
$ws1.Cells.Item(77, 5).Value = "LP1912"
$ws1.Cells.Item(78, 1).Value = "08:20:43"
$ws1.Cells.Item(78, 2).Value = "08:49"
$ws1.Cells.Item(78, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(78, 4).Value = 29
$ws1.Cells.Item(78, 5).Value = "LP1912"
$ws1.Cells.Item(79, 1).Value = "07:31:43"
$ws1.Cells.Item(79, 2).Value = "08:52"
$ws1.Cells.Item(79, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(79, 4).Value = 81
$ws1.Cells.Item(79, 5).Value = "LP1912"
$ws1.Cells.Item(80, 1).Value = "06:59:44"
$ws1.Cells.Item(80, 2).Value = "08:54"
$ws1.Cells.Item(80, 3).Value = "17_ROMERO"
$ws1.Cells.Item(80, 4).Value = 115
$ws1.Cells.Item(80, 5).Value = "LP1912"
$ws1.Cells.Item(81, 1).Value = "07:31:43"
$ws1.Cells.Item(81, 2).Value = "09:01"
$ws1.Cells.Item(81, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(81, 4).Value = 90
$ws1.Cells.Item(81, 5).Value = "LP1912"
$ws1.Cells.Item(82, 1).Value = "07:57:27"
$ws1.Cells.Item(82, 2).Value = "09:03"
$ws1.Cells.Item(82, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(82, 4).Value = 66
$ws1.Cells.Item(82, 5).Value = "LP1912"
$ws1.Cells.Item(83, 1).Value = "07:31:43"
$ws1.Cells.Item(83, 2).Value = "09:10"
$ws1.Cells.Item(83, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(83, 4).Value = 99
$ws1.Cells.Item(83, 5).Value = "LP1912"
$ws1.Cells.Item(84, 1).Value = "07:57:27"
$ws1.Cells.Item(84, 2).Value = "09:16"
$ws1.Cells.Item(84, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(84, 4).Value = 105
$ws1.Cells.Item(84, 5).Value = "LP1912"
$ws1.Cells.Item(85, 1).Value = "07:57:27"
$ws1.Cells.Item(85, 2).Value = "09:17"
$ws1.Cells.Item(85, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(85, 4).Value = 80
$ws1.Cells.Item(85, 5).Value = "LP1912"
$ws1.Cells.Item(86, 1).Value = "07:57:27"
$ws1.Cells.Item(86, 2).Value = "09:18"
$ws1.Cells.Item(86, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(86, 4).Value = 81
$ws1.Cells.Item(86, 5).Value = "LP1912"
$ws1.Cells.Item(87, 1).Value = "07:31:43"
$ws1.Cells.Item(87, 2).Value = "09:21"
$ws1.Cells.Item(87, 3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(87, 4).Value = 110
$ws1.Cells.Item(87, 5).Value = "LP1912"
$ws1.Cells.Item(88, 1).Value = "07:31:43"
$ws1.Cells.Item(88, 2).Value = "09:22"
$ws1.Cells.Item(88, 3).Value = "17_ROMERO"
$ws1.Cells.Item(88, 4).Value = 111
$ws1.Cells.Item(88, 5).Value = "LP1912"
$ws1.Cells.Item(89, 1).Value = "08:20:43"
$ws1.Cells.Item(89, 2).Value = "09:22"
$ws1.Cells.Item(89, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(89, 4).Value = 62
$ws1.Cells.Item(89, 5).Value = "LP1912"
$ws1.Cells.Item(90, 1).Value = "07:31:43"
$ws1.Cells.Item(90, 2).Value = "09:23"
$ws1.Cells.Item(90, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(90, 4).Value = 112
$ws1.Cells.Item(90, 5).Value = "LP1912"
$ws1.Cells.Item(91, 1).Value = "07:31:43"
$ws1.Cells.Item(91, 2).Value = "09:23"
$ws1.Cells.Item(91, 3).Value = "17_ROMERO"
$ws1.Cells.Item(91, 4).Value = 86
$ws1.Cells.Item(91, 5).Value = "LP1912"
$ws1.Cells.Item(92, 1).Value = "07:57:27"
$ws1.Cells.Item(92, 2).Value = "09:32"
$ws1.Cells.Item(92, 3).Value = "15_ABASTO"
$ws1.Cells.Item(92, 4).Value = 95
$ws1.Cells.Item(92, 5).Value = "LP1912"
$ws1.Cells.Item(93, 1).Value = "07:57:27"
$ws1.Cells.Item(93, 2).Value = "09:33"
$ws1.Cells.Item(93, 3).Value = "10_OLMOS"
$ws1.Cells.Item(93, 4).Value = 96
$ws1.Cells.Item(93, 5).Value = "LP1912"
$ws1.Cells.Item(94, 1).Value = "08:20:43"
$ws1.Cells.Item(94, 2).Value = "09:41"
$ws1.Cells.Item(94, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(94, 4).Value = 81
$ws1.Cells.Item(94, 5).Value = "LP1912"
$ws1.Cells.Item(95, 1).Value = "07:57:27"
$ws1.Cells.Item(95, 2).Value = "09:42"
$ws1.Cells.Item(95, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(95, 4).Value = 105
$ws1.Cells.Item(95, 5).Value = "LP1912"
$ws1.Cells.Item(96, 1).Value = "07:57:27"
$ws1.Cells.Item(96, 2).Value = "09:43"
$ws1.Cells.Item(96, 3).Value = "14_ABASTO"
$ws1.Cells.Item(96, 4).Value = 106
$ws1.Cells.Item(96, 5).Value = "LP1912"
$ws1.Cells.Item(97, 1).Value = "08:20:43"
$ws1.Cells.Item(97, 2).Value = "10:08"
$ws1.Cells.Item(97, 3).Value = "10_OLMOS"
$ws1.Cells.Item(97, 4).Value = 108
$ws1.Cells.Item(97, 5).Value = "LP1912"
$ws1.Cells.Item(98, 1).Value = "08:20:43"
$ws1.Cells.Item(98, 2).Value = "10:12"
$ws1.Cells.Item(98, 3).Value = "15_ABASTO"
$ws1.Cells.Item(98, 4).Value = 112
$ws1.Cells.Item(98, 5).Value = "LP1912"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = "Última actualización: 08:20:43"
$ws2.Cells.Item(3, 1).Value = "Total filas: 13"

$ws2.Cells.Item(14, 1).Value = "08:20:43"
$ws2.Cells.Item(14, 2).Value = "08:21"
$ws2.Cells.Item(14, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(14, 4).Value = 1
$ws2.Cells.Item(14, 5).Value = "LP1912"
$ws2.Cells.Item(15, 1).Value = "06:44:15"
$ws2.Cells.Item(15, 2).Value = "08:23"
$ws2.Cells.Item(15, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(15, 4).Value = 99
$ws2.Cells.Item(15, 5).Value = "LP1912"
$ws2.Cells.Item(16, 1).Value = "07:31:43"
$ws2.Cells.Item(16, 2).Value = "09:01"
$ws2.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(16, 4).Value = 90
$ws2.Cells.Item(16, 5).Value = "LP1912"
$ws2.Cells.Item(17, 1).Value = "08:20:43"
$ws2.Cells.Item(17, 2).Value = "09:41"
$ws2.Cells.Item(17, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(17, 4).Value = 81
$ws2.Cells.Item(17, 5).Value = "LP1912"
$ws2.Cells.Item(18, 1).Value = "07:57:27"
$ws2.Cells.Item(18, 2).Value = "09:42"
$ws2.Cells.Item(18, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(18, 4).Value = 105
$ws2.Cells.Item(18, 5).Value = "LP1912"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 1).Value = "Última actualización: 08:20:43"
$ws3.Cells.Item(3, 1).Value = "Total filas: 21"

$ws3.Cells.Item(23, 1).Value = "08:20:43"
$ws3.Cells.Item(23, 2).Value = "08:37"
$ws3.Cells.Item(23, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(23, 4).Value = 17
$ws3.Cells.Item(23, 5).Value = "L6173"
$ws3.Cells.Item(24, 1).Value = "07:31:43"
$ws3.Cells.Item(24, 2).Value = "09:08"
$ws3.Cells.Item(24, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(24, 4).Value = 97
$ws3.Cells.Item(24, 5).Value = "L6203"
$ws3.Cells.Item(25, 1).Value = "07:57:27"
$ws3.Cells.Item(25, 2).Value = "09:09"
$ws3.Cells.Item(25, 3).Value = "215D_LA PLATA"
$ws3.Cells.Item(25, 4).Value = 72
$ws3.Cells.Item(25, 5).Value = "L6203"
$ws3.Cells.Item(26, 1).Value = "08:20:43"
$ws3.Cells.Item(26, 2).Value = "10:02"
$ws3.Cells.Item(26, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(26, 4).Value = 102
$ws3.Cells.Item(26, 5).Value = "L6173"

